$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so numeric-looking
# values (e.g. "1.003") are not auto-converted into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.059.14'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '1.887.18'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").Value = '0.7377'
$ws.Range("E5").Value = '  -2.01%  '

$ws.Range("D6").Value = '242.82'
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("E7").Value = '  +0.33%  '

$ws.Range("D8").Value = '0.3166'
$ws.Range("E8").Value = '  +1.30%  '

$ws.Range("D9").Value = '0.07174'
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").Value = '24.71'
$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("D11").Value = '0.08349'
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").Value = '0.7566'
$ws.Range("E12").Value = '  -0.45%  '

$ws.Range("D13").Value = '5.415'
$ws.Range("E13").Value = '  +0.96%  '

$ws.Range("D14").Value = '1.812.31'
$ws.Range("E14").Value = '  -6.23%  '

$ws.Range("D15").Value = '92.87'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '6.151'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '30.062.37'
$ws.Range("E17").Value = '  +0.76%  '

$ws.Range("D18").Value = '250.28'
$ws.Range("E18").Value = '  +2.82%  '

$ws.Range("D19").Value = '13.58'
$ws.Range("E19").Value = '  -0.95%  '

$ws.Range("D20").Value = '0.000007862'
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Value = '2.191.60'
$ws.Range("E21").Value = '  +2.45%  '

$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '7.907'
$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("D25").Value = '0.1568'
$ws.Range("E25").Value = '  -1.73%  '

$ws.Range("D26").Value = '9.289'
$ws.Range("E26").Value = '  -0.91%  '

$ws.Range("D27").Value = '164.52'
$ws.Range("E27").Value = '  +0.99%  '

$ws.Range("D28").Value = '18.67'
$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("D29").Value = '2.048'
$ws.Range("E29").Value = '  +0.84%  '

$ws.Range("D30").Value = '1.478'
$ws.Range("E30").Value = '  -0.38%  '

$ws.Range("D31").Value = '4.565'
$ws.Range("E31").Value = '  +1.29%  '

$ws.Range("D32").Value = '1.537'
$ws.Range("E32").Value = '  +0.38%  '

$ws.Range("D33").Value = '4.194'
$ws.Range("E33").Value = '  +1.10%  '

$ws.Range("D34").Value = '0.05342'
$ws.Range("E34").Value = '  -1.60%  '

$ws.Range("D35").Value = '1.251'
$ws.Range("E35").Value = '  +0.71%  '

$ws.Range("D36").Value = '0.7678'
$ws.Range("E36").Value = '  +2.18%  '

$ws.Range("D37").Value = '1.002'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '2.732'
$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("D39").Value = '0.01957'
$ws.Range("E39").Value = '  +0.56%  '

$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").Value = '0.4554'
$ws.Range("E41").Value = '  +2.04%  '

$ws.Range("D42").Value = '1.102.30'
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("D43").Value = '6.058'
$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("D44").Value = '72.23'
$ws.Range("E44").Value = '  -0.76%  '

$ws.Range("D45").Value = '0.8742'
$ws.Range("E45").Value = '  +1.63%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.005'
$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '104.30'
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("D48").Value = '1.854'
$ws.Range("E48").Value = '  -0.31%  '

$ws.Range("D49").Value = '7.570'
$ws.Range("E49").Value = '  -1.93%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.589'
$ws.Range("E50").Value = '  -1.78%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.034.97'
$ws.Range("E51").Value = '  -0.08%  '
